# Insert a new accelerometer sample as row 2 (shifting the existing
# samples down by one row), and drop the final two samples so the
# sheet keeps the same number of rows (A1:C21 instead of A1:C22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2, pushing all existing data rows down.
$ws.Rows(2).Insert()

# Excel copies the formatting of the row above (the bold header) onto the
# newly inserted row; strip that back off so the new data row matches the
# plain formatting of the other data rows.
$ws.Rows(2).ClearFormats()

# Populate the new row 2 with the new x/y/z accelerometer reading.
$ws.Range("A2").Value = -3.373677730560303
$ws.Range("B2").Value = 8.367032051086426
$ws.Range("C2").Value = -1.811180233955384

# After the insert, the data that used to occupy the last two rows (old
# rows 21 and 22) now sits in rows 22 and 23. Remove those trailing rows
# so the sheet ends at row 21 again.
$ws.Range("A22:A23").EntireRow.Delete()
